# Added year to archetypes (only one construction per function for now)
# Fill in year_start / year_end columns (B, C) with "0" for every
# archetype row (2-19) on both the ARCHITECTURE and HVAC sheets.

$wb = $excel.ActiveWorkbook

$wsArchitecture = $wb.Worksheets.Item("ARCHITECTURE")
$wsHvac = $wb.Worksheets.Item("HVAC")

for ($row = 2; $row -le 19; $row++) {
    $wsArchitecture.Cells.Item($row, 2).Value = "0"
    $wsArchitecture.Cells.Item($row, 3).Value = "0"

    $wsHvac.Cells.Item($row, 2).Value = "0"
    $wsHvac.Cells.Item($row, 3).Value = "0"
}

# Match the cursor/selection state recorded in the target workbook:
# ARCHITECTURE ends up with C19 selected (and is no longer the active tab),
# while HVAC becomes the active tab with B19:C19 selected.
$wsArchitecture.Range("C19").Select()

$wsHvac.Activate()
$wsHvac.Range("B19:C19").Select()
